# Rename the "Slide_subtitle" placeholder text on the Title shape of slide 1
# to "Slide_title" (split across two runs: "Slide" + "_title"), matching the
# author's edit that re-typed part of the title's text.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Locate the title placeholder robustly (ppPlaceholderTitle = 1) instead of
# relying on a hard-coded shape index.
$titleShape = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $sh = $s.Shapes.Item($i)
    if ($sh.HasTextFrame -and $sh.Type -eq 14) {
        if ($sh.PlaceholderFormat.Type -eq 1) {
            $titleShape = $sh
            break
        }
    }
}
if ($titleShape -eq $null) {
    $titleShape = $s.Shapes.Item(2)
}

$tr = $titleShape.TextFrame.TextRange

# Split the run: keep "Slide" as-is, replace the remaining "_subtitle" with
# "_title" so the text reads "Slide_title" across two separate runs.
$head = $tr.Characters(1, 5)
$head.Text = "Slide"

$tail = $tr.Characters(6, $tr.Length - 5)
$tail.Text = "_title"
